# Weekly data refresh: a new Brócoli price record (most recent week) is
# inserted at the top of the data table (row 169), pushing all existing
# records for rows 169-203 down by one row (to 170-204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 169 - shifts rows 169:203 down to 170:204
# and extends the sheet's used range to A1:R204.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(169, 1).Value  = 7
$ws.Cells.Item(169, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value  = "Ñuble"
$ws.Cells.Item(169, 4).Value  = 44543
$ws.Cells.Item(169, 5).Value  = 16
$ws.Cells.Item(169, 6).Value  = 100112023
$ws.Cells.Item(169, 7).Value  = "Brócoli"
$ws.Cells.Item(169, 8).Value  = "Sin especificar"
$ws.Cells.Item(169, 9).Value  = "Primera"
$ws.Cells.Item(169, 10).Value = 400
$ws.Cells.Item(169, 11).Value = 700
$ws.Cells.Item(169, 12).Value = 800
$ws.Cells.Item(169, 13).Value = 750
$ws.Cells.Item(169, 14).Value = "$/unidad"
$ws.Cells.Item(169, 15).Value = "Región del Maule"
$ws.Cells.Item(169, 16).Value = 750
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"
